# Re-pulled dSF (column F) data: push all corrected values from the
# updated source, recalculating the mean-derived delta column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -1
    7  = -1
    14 = -2
    19 = 0
    20 = 2
    26 = -3
    27 = 1
    32 = 2
    33 = -5
    34 = -3
    35 = -7
    36 = -7
    39 = -3
    40 = -7
    41 = -6
    42 = -3
    43 = -4
    44 = -4
    46 = 3
    49 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
